$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.248.56'
$ws.Range("E2").Value = '  +0.28%  '
$ws.Range("D3").Value = '2.299.02'
$ws.Range("E3").Value = '  -0.18%  '
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '315.59'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.26%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '103.93'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.14%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.630'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.65%  '
$ws.Range("E8").Value = '  -0.08%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.608'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.13%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '39.84'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.97%  '
$ws.Range("E11").Value = '  -0.67%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '8.33'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.80%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.963'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.09%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '15.30'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.38%  '
$ws.Range("D16").Value = '2.650.37'
$ws.Range("E16").Value = '  +0.01%  '
$ws.Range("D17").Value = '2.307.18'
$ws.Range("E17").Value = '  +0.17%  '
$ws.Range("D18").Value = '42.358.04'
$ws.Range("E18").Value = '  +0.84%  '
$ws.Range("E19").Value = '  -2.17%  '
$ws.Range("E20").Value = '  +1.12%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '73.26'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.57%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '276.02'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +6.39%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.53'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.91%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.10'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +19.35%  '
$ws.Range("E25").Value = '  -1.12%  '
$ws.Range("E26").Value = '  -0.28%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.81'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.54%  '
$ws.Range("E28").Value = '  +3.32%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '22.76'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.03%  '
$ws.Range("E30").Value = '  +0.54%  '
$ws.Range("E31").Value = '  -0.16%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.0870'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.94%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.86'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.36%  '
$ws.Range("E34").Value = '  +3.91%  '
$ws.Range("E35").Value = '  -0.47%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.58'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -11.80%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0366'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +4.07%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.57'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.01%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.72'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +3.34%  '
$ws.Range("E40").Value = '  -0.99%  '
$ws.Range("E41").Value = '  +2.46%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '69.49'
$ws.Range("D42").Style = "Normal"
$ws.Range("B43").Value = 'BitcoinSV'
$ws.Range("C43").Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '94.93'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -3.12%  '
$ws.Range("B44").Value = 'Algorand'
$ws.Range("C44").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.226'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.82%  '
$ws.Range("E45").Value = '  +0.28%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '81.54'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +9.95%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '12.03'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.66%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '112.63'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.06%  '
$ws.Range("E50").Value = '  -2.67%  '
$ws.Range("D51").Value = '1.590.91'
$ws.Range("E51").Value = '  +1.86%  '
